$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.313941333333333
$ws.Range("H2").Value = 15.941824
$ws.Range("I2").Value = 0.176869630377001
$ws.Range("J2").Value = 0.176869630377001
$ws.Range("M2").Value = 0.9386610000000001
$ws.Range("N2").Value = 2.815983
$ws.Range("O2").Value = 0.04181245246793033
$ws.Range("P2").Value = 0.04181245246793032
$ws.Range("Q2").Value = 4.987989485888
$ws.Range("R2").Value = 44.891905372992
$ws.Range("S2").Value = 0.007395353013158758
$ws.Range("T2").Value = 0.007395353013158757
$ws.Range("G3").Value = 5.313941333333333
$ws.Range("H3").Value = 15.941824
$ws.Range("I3").Value = 0.176869630377001
$ws.Range("J3").Value = 0.176869630377001
$ws.Range("O3").Value = 0.1106393125456779
$ws.Range("P3").Value = 0.1106393125456779
$ws.Range("Q3").Value = 13.19864526308267
$ws.Range("R3").Value = 118.787807367744
$ws.Range("S3").Value = 0.01956873431511954
$ws.Range("T3").Value = 0.01956873431511954
$ws.Range("G4").Value = 5.313941333333333
$ws.Range("H4").Value = 15.941824
$ws.Range("I4").Value = 0.176869630377001
$ws.Range("J4").Value = 0.176869630377001
$ws.Range("O4").Value = 0.8475482349863918
$ws.Range("P4").Value = 0.8475482349863918
$ws.Range("Q4").Value = 101.107718762432
$ws.Range("R4").Value = 909.9694688618879
$ws.Range("S4").Value = 0.1499055430487227
$ws.Range("T4").Value = 0.1499055430487227
$ws.Range("I5").Value = 0.5461014638447835
$ws.Range("J5").Value = 0.5461014638447835
$ws.Range("M5").Value = 0.9386610000000001
$ws.Range("N5").Value = 2.815983
$ws.Range("O5").Value = 0.04181245246793033
$ws.Range("P5").Value = 0.04181245246793032
$ws.Range("Q5").Value = 15.400882300029
$ws.Range("R5").Value = 138.607940700261
$ws.Range("S5").Value = 0.02283384149967718
$ws.Range("T5").Value = 0.02283384149967718
$ws.Range("I6").Value = 0.5461014638447835
$ws.Range("J6").Value = 0.5461014638447835
$ws.Range("O6").Value = 0.1106393125456779
$ws.Range("P6").Value = 0.1106393125456779
$ws.Range("S6").Value = 0.06042029053997523
$ws.Range("T6").Value = 0.06042029053997522
$ws.Range("I7").Value = 0.5461014638447835
$ws.Range("J7").Value = 0.5461014638447835
$ws.Range("O7").Value = 0.8475482349863918
$ws.Range("P7").Value = 0.8475482349863918
$ws.Range("S7").Value = 0.4628473318051311
$ws.Range("T7").Value = 0.4628473318051311
$ws.Range("I8").Value = 0.2770289057782155
$ws.Range("J8").Value = 0.2770289057782155
$ws.Range("M8").Value = 0.9386610000000001
$ws.Range("N8").Value = 2.815983
$ws.Range("O8").Value = 0.04181245246793033
$ws.Range("P8").Value = 0.04181245246793032
$ws.Range("Q8").Value = 7.812631633613
$ws.Range("R8").Value = 70.31368470251699
$ws.Range("S8").Value = 0.01158325795509439
$ws.Range("T8").Value = 0.01158325795509439
$ws.Range("I9").Value = 0.2770289057782155
$ws.Range("J9").Value = 0.2770289057782155
$ws.Range("O9").Value = 0.1106393125456779
$ws.Range("P9").Value = 0.1106393125456779
$ws.Range("S9").Value = 0.03065028769058314
$ws.Range("T9").Value = 0.03065028769058314
$ws.Range("I10").Value = 0.2770289057782155
$ws.Range("J10").Value = 0.2770289057782155
$ws.Range("O10").Value = 0.8475482349863918
$ws.Range("P10").Value = 0.8475482349863918
$ws.Range("S10").Value = 0.234795360132538
$ws.Range("T10").Value = 0.234795360132538
